$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.289.65"
$ws.Range("E2").Value = "  -3.04%  "
$ws.Range("D3").Value = "1.771.65"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "'305.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("D7").Value = "'0.4229"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.21%  "
$ws.Range("D8").Value = "'0.3603"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.41%  "
$ws.Range("D9").Value = "'0.07123"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").Value = "'0.8365"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("D11").Value = "'20.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "1.765.13"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("D13").Value = "'6.442"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("D14").Value = "'5.238"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("D15").Value = "'0.06888"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("D18").Value = "'0.000008631"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").Value = "'1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "'14.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("D21").Value = "26.296.73"
$ws.Range("E21").Value = "  -2.85%  "
$ws.Range("D22").Value = "'5.101"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("D23").Value = "'10.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").Value = "1.975.22"
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").Value = "'151.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("D26").Value = "'1.795"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.90%  "
$ws.Range("D27").Value = "'18.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").Value = "'5.056"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("D29").Value = "'114.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").Value = "'1.834"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.07%  "
$ws.Range("D31").Value = "'0.08815"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").Value = "'0.7245"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").Value = "'1.115"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.46%  "
$ws.Range("D34").Value = "'4.311"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").Value = "'1.000"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("D36").Value = "'2.734"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.10%  "
$ws.Range("D37").Value = "'1.103"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.85%  "
$ws.Range("D38").Value = "'0.05103"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").Value = "'0.01884"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("D40").Value = "'0.1609"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("D41").Value = "'0.4914"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("D42").Value = "'2.592"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "'6.342"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.78%  "
$ws.Range("B44").Value = "PaxosStandard"
$ws.Range("C44").Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
$ws.Range("D44").Value = "'1.003"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -30.96%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'8.024"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'104.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'10.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.616"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06168"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.4439"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.36%  "
